{"js": "// Locate the paragraph that starts with the \"Atendiendo a lo establecido...\"\n// data-protection notice, trim the stray trailing \"s.\" so it reads\n// \"...y sus Municipio\" instead of \"...y sus Municipios.\", and remove the\n// two blank paragraphs that immediately follow it (they were left over\n// spacer paragraphs that caused extra blank space / pages when the\n// template body was concatenated with other generated documents).\nconst marker = \"Atendiendo a lo establecido\";\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nlet target = null;\nlet targetIndex = -1;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.indexOf(marker) === 0) {\n    target = paragraphs.items[i];\n    targetIndex = i;\n    break;\n  }\n}\n\nif (target) {\n  const oldText = target.text;\n  const newText = oldText.slice(0, -2); // drop the trailing \"s.\"\n  target.insertText(newText, Word.InsertLocation.replace);\n  await context.sync();\n\n  // Re-fetch the paragraph collection since the one we just edited was\n  // replaced; grab the two paragraphs right after it (both empty\n  // spacer paragraphs) and delete them.\n  const refreshed = context.document.body.paragraphs;\n  refreshed.load(\"text\");\n  await context.sync();\n\n  const firstEmpty = refreshed.items[targetIndex + 1];\n  const secondEmpty = refreshed.items[targetIndex + 1];\n  // Delete the paragraph right after the target twice: once it is\n  // removed, the following paragraph shifts into that same slot.\n  if (firstEmpty) {\n    firstEmpty.delete();\n    await context.sync();\n  }\n  const refreshed2 = context.document.body.paragraphs;\n  refreshed2.load(\"text\");\n  await context.sync();\n  const nextEmpty = refreshed2.items[targetIndex + 1];\n  if (nextEmpty) {\n    nextEmpty.delete();\n    await context.sync();\n  }\n}\n", "ps1": "# Locate the paragraph that starts with the \"Atendiendo a lo establecido...\"\n# data-protection notice, trim the stray trailing \"s.\" so it reads\n# \"...y sus Municipio\" instead of \"...y sus Municipios.\", and remove the\n# two blank paragraphs that immediately follow it (left-over spacer\n# paragraphs that produced extra blank space / pages when this template's\n# body was concatenated with other generated documents).\n$d = $word.ActiveDocument\n\n$marker = \"Atendiendo a lo establecido\"\n$targetIndex = -1\n$count = $d.Paragraphs.Count\nfor ($i = 1; $i -le $count; $i++) {\n    if ($d.Paragraphs.Item($i).Range.Text.StartsWith($marker)) {\n        $targetIndex = $i\n        break\n    }\n}\n\nif ($targetIndex -ne -1) {\n    $targetRange = $d.Paragraphs.Item($targetIndex).Range\n\n    $find = $targetRange.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = \"y sus Municipios.\"\n    $find.Replacement.Text = \"y sus Municipio\"\n    $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n\n    # Remove the two empty spacer paragraphs right after the target paragraph.\n    $d.Paragraphs.Item($targetIndex + 1).Range.Delete()\n    $d.Paragraphs.Item($targetIndex + 1).Range.Delete()\n}\n"}
